$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------
# The sheet is a single column (A) list of ISBNs. Four new ISBN rows
# were inserted (pushing everything below them down), one of which was
# typed as a bare number (kept numeric) while the other three (plus one
# pre-existing "text-look-alike" 12-digit code) were entered/forced as
# Text-formatted cells.
#
# Final row layout (1-based):
#   10 -> "979-1-876543-21-0"      (new, default format)
#   12 -> 5556667778               (new, Text format, stored as number)
#   17 -> "9780307455376"          (new, Text format)
#   25 -> "9788423687268"          (new, Text format)
# everything else shifts down accordingly.
#
# Insert rows top-to-bottom at their FINAL target row numbers so each
# insert only pushes the still-unplaced rows further down.
# --------------------------------------------------------------------

$ws.Rows.Item(10).EntireRow.Insert()
$ws.Rows.Item(12).EntireRow.Insert()
$ws.Rows.Item(17).EntireRow.Insert()
$ws.Rows.Item(25).EntireRow.Insert()

# Row 17 - new ISBN, forced to Text format before typing so it lands in
# the shared-string table rather than being parsed as a number.
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "9780307455376"

# Row 10 - new ISBN with dashes; stays text naturally.
$ws.Range("A10").Value = "979-1-876543-21-0"

# Row 25 - new ISBN, forced to Text format like row 17.
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "9788423687268"

# Row 12 - typed as a genuine number, then the cell format switched to
# Text afterwards (value stays numeric, only the display format changes).
$ws.Range("A12").Value = 5556667778
$ws.Range("A12").NumberFormat = "@"

# Column A was re-fit to the new (slightly narrower) content
# (17.28515625 -> 16.85546875 character-width units).
$ws.Columns.Item(1).ColumnWidth = 16

# Leave the cursor on the last-touched cell, matching the editor's final
# active cell in the source session.
$ws.Range("D8").Select()
